# Add exportDateFormat() option for XLS/XLSX files to ensure dates are
# formatted correctly for PHP.
#
# Adds a "date" column (header + two date values) to Sheet1 column E,
# makes Sheet1 the active sheet with E1:E3 selected, and clears the
# previous "Third Sheet" tab selection.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Header cell E1: "date", bold ---
$rngE1 = $ws1.Range("E1")
$rngE1.ClearFormats()
$rngE1.Value = "date"
$rngE1.Font.Bold = $true

# --- Data cell E2: date value formatted as a short date ---
$rngE2 = $ws1.Range("E2")
$rngE2.ClearFormats()
$rngE2.Value = 44044
$rngE2.NumberFormat = "mm-dd-yy"

# --- Data cell E3: same date value / format as E2 (copy format so both
#     cells share the exact same style record) ---
$rngE3 = $ws1.Range("E3")
$rngE3.Value = 44044
$rngE2.Copy() | Out-Null
$rngE3.PasteSpecial(-4122) | Out-Null

# --- Make Sheet1 the active sheet with E1:E3 selected (active cell E1) ---
$ws1.Activate() | Out-Null
$ws1.Range("E1:E3").Select() | Out-Null
